$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 15
$ws.Range("C1").Value = 16
$ws.Range("D1").Value = 15
$ws.Range("E1").Value = 16

$ws.Range("B2").Value = 289.32168099038631
$ws.Range("C2").Value = 261.58384151227301
$ws.Range("D2").Value = 289.58720152724482
$ws.Range("E2").Value = 257.34769713850034

$ws.Range("B3").Value = 294.73978121257522
$ws.Range("C3").Value = 250.20360560472221
$ws.Range("D3").Value = 305.39943412595017
$ws.Range("E3").Value = 248.23193479038781

$ws.Range("B1:E3").Select()
